$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 9 (Политика №8): I9/J9 go from 2/2 to 0/0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0

# Fix existing row 10 (Политика №9): I10/J10 go from 4/4 to 0/0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0

# New template rows 11-16
$newRows = @(
    @{ Row = 11; Policy = "Политика №10"; Task = "Задание №10"; E = 0; F = 0; G = 0; H = 0; I = 0; J = 0 },
    @{ Row = 12; Policy = "Политика №11"; Task = "Задание №11"; E = 0; F = 0; G = 0; H = 0; I = 0; J = 0 },
    @{ Row = 13; Policy = "Политика №12"; Task = "Задание №12"; E = 1; F = 0; G = 4; H = 0; I = 0; J = 5 },
    @{ Row = 14; Policy = "Политика №13"; Task = "Задание №13"; E = 1; F = 0; G = 1; H = 0; I = 0; J = 2 },
    @{ Row = 15; Policy = "Политика №14"; Task = "Задание №14"; E = 3; F = 1; G = 2; H = 1; I = 2; J = 9 },
    @{ Row = 16; Policy = "Политика №15"; Task = "Задание №15"; E = 0; F = 0; G = 0; H = 0; I = 0; J = 0 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("C$row").Value = $r.Policy
    $ws.Range("D$row").Value = $r.Task
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J

    # Match the bordered style used by the existing data rows
    $ws.Range("C$row`:J$row").Borders.LineStyle = 1
}
